$p = $ppt.ActivePresentation

# =========================================================================
# Slide 1: title "Sheet1" -> "Full portfoilio Allocation"; table expanded
# from a single "Asset" column into 4 columns (Asset/Sector/Country/Value)
# and the leading "2. Portfolio Allocation Data" banner row is dropped.
# =========================================================================
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Full portfoilio Allocation"

$tbl1 = $s1.Shapes.Item(2).Table

# Drop the old "2. Portfolio Allocation Data" banner row.
$tbl1.Rows.Item(1).Delete()

# Grow the single "Asset" column into four columns.
$tbl1.Columns.Add() | Out-Null
$tbl1.Columns.Add() | Out-Null
$tbl1.Columns.Add() | Out-Null

$data1 = @(
    @("Asset",   "Sector", "Country", "Value"),
    @("AAPL",    "Tech",   "USA",     "20000"),
    @("TSLA",    "Auto",   "USA",     "15000"),
    @("TCS.NS",  "IT",     "India",   "10000"),
    @("BTC-USD", "Crypto", "Global",  "8000"),
    @("ETH-USD", "Crypto", "Global",  "5000")
)

for ($r = 1; $r -le $tbl1.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl1.Columns.Count; $c++) {
        $tbl1.Cell($r, $c).Shape.TextFrame.TextRange.Text = $data1[$r - 1][$c - 1]
    }
    $tbl1.Rows.Item($r).Height = 48
}

for ($c = 1; $c -le $tbl1.Columns.Count; $c++) {
    $tbl1.Columns.Item($c).Width = 162
}

# =========================================================================
# Slide 2: title "Full portfoilio Allocation" -> " Top 3 Portfoilio Assests";
# table trimmed down to the top 3 assets and expanded to 4 columns.
# =========================================================================
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = " Top 3 Portfoilio Assests"

$tbl2 = $s2.Shapes.Item(2).Table

# Drop the old "2. Portfolio Allocation Data" banner row.
$tbl2.Rows.Item(1).Delete()
# Keep only the top 3 assets (AAPL, TSLA, TCS.NS) - drop BTC-USD / ETH-USD.
$tbl2.Rows.Item(6).Delete()
$tbl2.Rows.Item(5).Delete()

$tbl2.Columns.Add() | Out-Null
$tbl2.Columns.Add() | Out-Null
$tbl2.Columns.Add() | Out-Null

$data2 = @(
    @("Asset",  "Sector", "Country", "Value"),
    @("AAPL",   "Tech",   "USA",     "20000"),
    @("TSLA",   "Auto",   "USA",     "15000"),
    @("TCS.NS", "IT",     "India",   "10000")
)

for ($r = 1; $r -le $tbl2.Rows.Count; $r++) {
    for ($c = 1; $c -le $tbl2.Columns.Count; $c++) {
        $tbl2.Cell($r, $c).Shape.TextFrame.TextRange.Text = $data2[$r - 1][$c - 1]
    }
    $tbl2.Rows.Item($r).Height = 72
}

for ($c = 1; $c -le $tbl2.Columns.Count; $c++) {
    $tbl2.Columns.Item($c).Width = 162
}

# =========================================================================
# Slide 3 ("Top 3 Portfoilio Assests" stub) is no longer needed now that
# its content lives on slide 2 - remove it from the deck entirely.
# =========================================================================
$p.Slides.Item(3).Delete()
